$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("G3").Value = 2.55
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.7
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 3.25
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.67
$ws.Range("V3").Value = 2.1
$ws.Range("X3").Value = 13
$ws.Range("Y3").Value = 10
$ws.Range("AC3").Value = 11
$ws.Range("AH3").Value = 10
$ws.Range("AJ3").Value = 10
$ws.Range("AK3").Value = 26
$ws.Range("AP3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 151
$ws.Range("AT3").Value = 2.75
$ws.Range("AW3").Value = 4.75

# Row 4 updates
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.6

# Row 5 updates
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
